# Add a new forecast-date column (Z) and a new observation row (38) to both
# the "cases" and "deaths" sheets, and backfill the previously-missing
# B24 "Observed" value on each sheet.
#
# Note: the values in column A (row labels) and row 1 (column headers) are
# stored as plain TEXT ("2020-05-19"), not real Excel dates. Assigning a
# date-shaped string straight to .Value would make Excel auto-convert it to
# a date serial number, so every such text assignment below first forces the
# cell to the "@" (Text) number format, then restores the "Normal" style
# once the text value has been committed (this keeps the cell itself free
# of any leftover explicit style index, matching the rest of the sheet).

$wb = $excel.ActiveWorkbook

$newDateLabel = "2020-05-19"
$newRow = 38
$newCol = 26   # column Z

# Per-sheet data: the new Z-column values for rows 25-38 (the "staircase"
# of forecasts as of the new date), plus the backfilled B24 value.
$casesZ = @{
    25 = 13356; 26 = 14220; 27 = 15253; 28 = 15996; 29 = 17129; 30 = 17940;
    31 = 18804; 32 = 19647; 33 = 20455; 34 = 20958; 35 = 21701; 36 = 22226;
    37 = 22693; 38 = 23219
}
$casesB24 = 12391

$deathsZ = @{
    25 = 1186; 26 = 1250; 27 = 1305; 28 = 1374; 29 = 1415; 30 = 1477;
    31 = 1563; 32 = 1626; 33 = 1685; 34 = 1718; 35 = 1767; 36 = 1803;
    37 = 1855; 38 = 1889
}
$deathsB24 = 1123

$sheetSpecs = @(
    @{ Name = "cases";  Z = $casesZ;  B24 = $casesB24 },
    @{ Name = "deaths"; Z = $deathsZ; B24 = $deathsB24 }
)

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

foreach ($spec in $sheetSpecs) {
    $ws = $wb.Worksheets.Item($spec.Name)

    # New column header (row 1) -> the new forecast "as of" date (text).
    Set-TextValue $ws.Range("Z1") $newDateLabel

    # New (empty) column-Z cells for the existing observed rows 2-24 -- no
    # forecast data exists yet for this column at those observation dates.
    for ($r = 2; $r -le 24; $r++) {
        $cell = $ws.Cells.Item($r, $newCol)
        $cell.NumberFormat = "@"
        $cell.Style = "Normal"
    }

    # New column-Z forecast data, rows 25-38 (the staircase of values).
    foreach ($r in ($spec.Z.Keys | Sort-Object)) {
        $ws.Cells.Item($r, $newCol).Value = $spec.Z[$r]
    }

    # Backfill the previously blank "Observed" value for 2020-05-05 (B24).
    $ws.Range("B24").Value = $spec.B24

    # New row 38: the observed-date label in column A (text)...
    Set-TextValue $ws.Range("A38") $newDateLabel

    # ...empty cells for columns B through Y (no data yet for this new
    # observed date in those older forecast columns)...
    for ($c = 2; $c -le 25; $c++) {
        $cell = $ws.Cells.Item($newRow, $c)
        $cell.NumberFormat = "@"
        $cell.Style = "Normal"
    }

    # ...and the single new-column (Z) forecast value for the new row.
    $ws.Cells.Item($newRow, $newCol).Value = $spec.Z[$newRow]
}
